# Meilenstein 3 - "Eine Folie der PPP ergänzt"
# On slide 7 ("Technologien"), the content placeholder ends with a
# paragraph that only contains a tab character. The author added a new
# bullet line reading "Log4j" by typing it in front of that trailing tab,
# which splits that paragraph into two runs: "Log4j" followed by the
# pre-existing tab run.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

$shp = $null
foreach ($candidate in $s.Shapes) {
    if ($candidate.Name -eq "Inhaltsplatzhalter 2") {
        $shp = $candidate
        break
    }
}

$tr = $shp.TextFrame.TextRange
$lastPara = $tr.Paragraphs($tr.Paragraphs().Count, 1)
[void]$lastPara.InsertBefore("Log4j")
